$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Marque") to make room for "Modele"
$ws.Columns.Item(2).Insert()

# The old Stock column (F) was pushed to G by the insert; remove the
# now-duplicate trailing column.
$ws.Columns.Item(7).Delete()

# Set header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Modele"
$ws.Range("C1").Value = "Marque"
$ws.Range("D1").Value = "Adresse Ip"
$ws.Range("E1").Value = "Departement"
$ws.Range("F1").Value = "Stock"

# Row 2
$ws.Range("A2").Value = 47
$ws.Range("B2").Value = "Epson EcoTank"
$ws.Range("C2").Value = "hp "
$ws.Range("D2").Value = "168.192.102.28"
$ws.Range("E2").Value = "exports"
$ws.Range("F2").Value = 3

# Row 3
$ws.Range("A3").Value = 48
$ws.Range("B3").Value = "Epson laser 18"
$ws.Range("C3").Value = "hp"
$ws.Range("D3").Value = "168.12.1.10"
$ws.Range("E3").Value = "achates"
$ws.Range("F3").Value = 6

# Row 4
$ws.Range("A4").Value = 85
$ws.Range("B4").Value = "canon i-sensys lbp12x"
$ws.Range("C4").Value = "canon"
$ws.Range("D4").Value = "168.12.10.5"
$ws.Range("E4").Value = "productione"
$ws.Range("F4").Value = 7

# Remove old rows 5 and 6 (previously held data that is no longer present)
$ws.Range("A5:F6").Delete()

$ws.Range("F4").Select()
